$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.337.67'
$ws.Range("E2").Value = '  +1.03%  '

# Row 3
$ws.Range("D3").Value = '1.666.80'
$ws.Range("E3").Value = '  +0.98%  '

# Row 4
$ws.Range("E4").Value = '  +0.91%  '

# Row 5
$ws.Range("D5").Value = '''219.30'
$ws.Range("E5").Value = '  +0.89%  '

# Row 6
$ws.Range("D6").Value = '''0.5347'
$ws.Range("E6").Value = '  +1.47%  '

# Row 7
$ws.Range("E7").Value = '  +0.86%  '

# Row 8
$ws.Range("D8").Value = '''0.2662'
$ws.Range("E8").Value = '  +2.64%  '

# Row 9
$ws.Range("D9").Value = '''0.06388'
$ws.Range("E9").Value = '  +1.28%  '

# Row 10
$ws.Range("D10").Value = '''20.85'
$ws.Range("E10").Value = '  +2.68%  '

# Row 11
$ws.Range("D11").Value = '''0.07845'
$ws.Range("E11").Value = '  +0.65%  '

# Row 12
$ws.Range("D12").Value = '''4.554'
$ws.Range("E12").Value = '  +0.88%  '

# Row 13
$ws.Range("D13").Value = '1.669.64'
$ws.Range("E13").Value = '  +1.38%  '

# Row 14
$ws.Range("D14").Value = '1.894.65'
$ws.Range("E14").Value = '  +0.95%  '

# Row 15
$ws.Range("D15").Value = '''0.5543'
$ws.Range("E15").Value = '  +1.24%  '

# Row 16
$ws.Range("D16").Value = '0.0₅8189'
$ws.Range("E16").Value = '  +0.12%  '

# Row 17
$ws.Range("D17").Value = '''65.92'
$ws.Range("E17").Value = '  +0.74%  '

# Row 18
$ws.Range("D18").Value = '26.344.93'
$ws.Range("E18").Value = '  +1.08%  '

# Row 19
$ws.Range("E19").Value = '  +0.89%  '

# Row 20
$ws.Range("E20").Value = '  +2.25%  '

# Row 21
$ws.Range("D21").Value = '''193.89'
$ws.Range("E21").Value = '  +1.83%  '

# Row 22
$ws.Range("D22").Value = '''10.28'
$ws.Range("E22").Value = '  +2.04%  '

# Row 23
$ws.Range("D23").Value = '''6.043'
$ws.Range("E23").Value = '  +0.51%  '

# Row 24
$ws.Range("E24").Value = '  +0.87%  '

# Row 25
$ws.Range("D25").Value = '''146.04'
$ws.Range("E25").Value = '  +1.78%  '

# Row 26
$ws.Range("E26").Value = '  -0.55%  '

# Row 27
$ws.Range("D27").Value = '''7.211'
$ws.Range("E27").Value = '  -0.04%  '

# Row 28
$ws.Range("D28").Value = '''16.17'
$ws.Range("E28").Value = '  +1.20%  '

# Row 29
$ws.Range("E29").Value = '  +4.46%  '

# Row 30
$ws.Range("D30").Value = '''0.05866'
$ws.Range("E30").Value = '  +1.14%  '

# Row 31
$ws.Range("D31").Value = '''1.283'
$ws.Range("E31").Value = '  +1.03%  '

# Row 32
$ws.Range("D32").Value = '''3.592'
$ws.Range("E32").Value = '  +1.36%  '

# Row 33
$ws.Range("D33").Value = '''3.287'

# Row 34
$ws.Range("D34").Value = '''1.609'
$ws.Range("E34").Value = '  +1.13%  '

# Row 35
$ws.Range("D35").Value = '''0.9695'
$ws.Range("E35").Value = '  +2.98%  '

# Row 36
$ws.Range("D36").Value = '''2.831'
$ws.Range("E36").Value = '  +1.34%  '

# Row 37
$ws.Range("E37").Value = '  +0.30%  '

# Row 38
$ws.Range("D38").Value = '''0.5821'
$ws.Range("E38").Value = '  +1.40%  '

# Row 40
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '1.068.13'
$ws.Range("E40").Value = '  +3.92%  '

# Row 41
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '''0.8621'
$ws.Range("E41").Value = '  +1.55%  '

# Row 42
$ws.Range("D42").Value = '''5.839'
$ws.Range("E42").Value = '  +2.29%  '

# Row 43
$ws.Range("E43").Value = '  +0.90%  '

# Row 44
$ws.Range("D44").Value = '''104.26'
$ws.Range("E44").Value = '  -0.85%  '

# Row 45
$ws.Range("D45").Value = '1.804.76'
$ws.Range("E45").Value = '  +0.68%  '

# Row 46
$ws.Range("D46").Value = '''58.07'
$ws.Range("E46").Value = '  +1.71%  '

# Row 47
$ws.Range("E47").Value = '  +1.59%  '

# Row 48
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₈104'
$ws.Range("E48").Value = '  -6.05%  '

# Row 49
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").Value = '''0.4393'
$ws.Range("E49").Value = '  +1.46%  '

# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '''8.005'
$ws.Range("E50").Value = '  +2.38%  '

# Row 51
$ws.Range("D51").Value = '''0.05165'
$ws.Range("E51").Value = '  +0.44%  '
